$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '61.683.52'
$ws.Range("E2").Value = '  +1.31%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.396.14'
$ws.Range("E3").Value = '  +0.11%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.998'
$ws.Range("E4").Value = '  -0.13%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '577.27'
$ws.Range("E5").Value = '  +1.28%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '140.72'
$ws.Range("E6").Value = '  -0.64%  '

$ws.Range("E7").Value = '  -0.03%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.475'
$ws.Range("E8").Value = '  +0.10%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.71'
$ws.Range("E9").Value = '  +2.66%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.123'
$ws.Range("E10").Value = '  -0.76%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.387'
$ws.Range("E11").Value = '  -1.91%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '3.967.34'
$ws.Range("E12").Value = '  -0.05%  '

$ws.Range("B13").Value = 'Avalanche'
$ws.Range("C13").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '28.74'
$ws.Range("E13").Value = '  +0.85%  '

$ws.Range("B14").Value = 'TRON'
$ws.Range("C14").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.126'
$ws.Range("E14").Value = '  +0.76%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.379.06'
$ws.Range("E15").Value = '  -0.21%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000170'
$ws.Range("E16").Value = '  -0.37%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '61.732.21'
$ws.Range("E17").Value = '  +1.28%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.14'
$ws.Range("E18").Value = '  -1.23%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.66'
$ws.Range("E19").Value = '  -2.33%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.97'
$ws.Range("E20").Value = '  -0.54%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '390.56'
$ws.Range("E21").Value = '  +1.75%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '75.43'
$ws.Range("E22").Value = '  +2.02%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.557'
$ws.Range("E23").Value = '  -0.52%  '

$ws.Range("E24").Value = '  -0.07%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000113'
$ws.Range("E25").Value = '  -3.44%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.190'
$ws.Range("E26").Value = '  +6.50%  '

$ws.Range("E27").Value = '  +0.01%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.28'
$ws.Range("E28").Value = '  -1.76%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.07'
$ws.Range("E29").Value = '  +1.00%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.14'
$ws.Range("E30").Value = '  -0.27%  '

$ws.Range("E31").Value = '  -0.06%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.37'
$ws.Range("E32").Value = '  -3.84%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '23.39'
$ws.Range("E33").Value = '  -1.05%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.93'
$ws.Range("E34").Value = '  -1.00%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '167.61'
$ws.Range("E35").Value = '  +1.27%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.03'
$ws.Range("E36").Value = '  +0.92%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.426.36'
$ws.Range("E37").Value = '  +0.15%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.47'
$ws.Range("E38").Value = '  -0.79%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0767'
$ws.Range("E39").Value = '  -0.82%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '25.98'
$ws.Range("E40").Value = '  -7.45%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.779'
$ws.Range("E41").Value = '  -0.01%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.42'
$ws.Range("E42").Value = '  -0.11%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.67'
$ws.Range("E43").Value = '  +0.04%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.13'
$ws.Range("E44").Value = '  +0.33%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.468.59'
$ws.Range("E45").Value = '  -0.80%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '6.70'
$ws.Range("E46").Value = '  -1.70%  '

$ws.Range("B47").Value = 'InjectiveProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '22.48'
$ws.Range("E47").Value = '  -4.08%  '

$ws.Range("B48").Value = 'FirstDigitalUSD'
$ws.Range("C48").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.996'
$ws.Range("E48").Value = '  -0.33%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0261'
$ws.Range("E49").Value = '  -3.78%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.03'
$ws.Range("E50").Value = '  -2.32%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.206'
$ws.Range("E51").Value = '  -1.43%  '
